$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45188 (2023-09-19) to 45189 (2023-09-20)
# for every data row (rows 2 through 408).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 408) { $lastRow = 408 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45189
